$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 <- values that used to belong to row 16 (A,B,D,E,F,G,H,Q,R)
$ws.Range("A14").Value = 130930230
$ws.Range("B14").Value = 78255
$ws.Range("E14").Value = 228579
$ws.Range("F14").Value = "Liten svartspik"
$ws.Range("G14").Value = "Chaenothecopsis nana"
$ws.Range("H14").Value = "Tibell"
$ws.Range("Q14").Value = 448404
$ws.Range("R14").Value = 7037411

# Row 15 <- values that used to belong to row 14
$ws.Range("A15").Value = 130930220
$ws.Range("B15").Value = 79714
$ws.Range("D15").Value = "NT"
$ws.Range("E15").Value = 1797
$ws.Range("F15").Value = "Mjölig dropplav"
$ws.Range("G15").Value = "Cliostomum leprosum"
$ws.Range("H15").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("Q15").Value = 448353
$ws.Range("R15").Value = 7037267

# Row 16 <- values that used to belong to row 15
$ws.Range("A16").Value = 130930219
$ws.Range("B16").Value = 92530
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 3298
$ws.Range("F16").Value = "Trådticka"
$ws.Range("G16").Value = "Climacocystis borealis"
$ws.Range("H16").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q16").Value = 448355
$ws.Range("R16").Value = 7037273

# Row 20 <-> Row 21 swap (A,B,D,E,F,G,H,Q,R)
$ws.Range("A20").Value = 130930222
$ws.Range("B20").Value = 83221
$ws.Range("D20").Value = "VU"
$ws.Range("E20").Value = 6486
$ws.Range("F20").Value = "Skuggnål"
$ws.Range("G20").Value = "Chaenotheca sphaerocephala"
$ws.Range("H20").Value = "Nádv."
$ws.Range("Q20").Value = 448330
$ws.Range("R20").Value = 7037323

$ws.Range("A21").Value = 130930223
$ws.Range("B21").Value = 79714
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 1797
$ws.Range("F21").Value = "Mjölig dropplav"
$ws.Range("G21").Value = "Cliostomum leprosum"
$ws.Range("H21").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("Q21").Value = 448337
$ws.Range("R21").Value = 7037328

# Row 25 <-> Row 26 swap (A, Q, R, AC)
$ws.Range("A25").Value = 131064804
$ws.Range("Q25").Value = 448308
$ws.Range("R25").Value = 7037158
$ws.Range("AC25").Value = "Ringhack äldre"

$ws.Range("A26").Value = 131064799
$ws.Range("Q26").Value = 448242
$ws.Range("R26").Value = 7037242
$ws.Range("AC26").Value = "Ringhack färska och äldre"
